$wb = $excel.ActiveWorkbook

# --- "Question Validation Succeed" sheet: add two new rows of test data ---
$succeed = $wb.Worksheets.Item("Question Validation Succeed")

# Write the new unique shared strings in the exact order needed so they land
# at the same shared-string-table indices as the target workbook:
#   141 Succeed-validationCriteria-mandatory-2
#   142 Succeed-validationCriteria-mandatory-1
#   143 { "mandatory": {"encounterType":"admission"} }
#   144 { "mandatory": {"encounterType":["admission","surveyResponse"]} }
$succeed.Range("A12").Value = "Succeed-validationCriteria-mandatory-2"
$succeed.Range("A11").Value = "Succeed-validationCriteria-mandatory-1"
$succeed.Range("K11").Value = '{ "mandatory": {"encounterType":"admission"} }'
$succeed.Range("K12").Value = '{ "mandatory": {"encounterType":["admission","surveyResponse"]} }'

# Fill in the remaining columns for both new rows (reuse existing strings).
$succeed.Range("B11").Value = "SurveyAnswer"
$succeed.Range("C11").Value = "SurveyAnswer: Full config"
$succeed.Range("P11").Value = '{ "source": "xyz" }'

$succeed.Range("B12").Value = "SurveyAnswer"
$succeed.Range("C12").Value = "SurveyAnswer: Full config"
$succeed.Range("P12").Value = '{ "source": "xyz" }'

# Widen column A to fit the new, longer row labels.
$succeed.Columns.Item(1).ColumnWidth = 56

# Configure the page setup for this sheet.
$succeed.PageSetup.PaperSize = 9
$succeed.PageSetup.Orientation = 1

# --- "Question Validation Fail" sheet: move the saved selection ---
$fail = $wb.Worksheets.Item("Question Validation Fail")
$fail.Activate() | Out-Null
$fail.Range("K12").Select() | Out-Null

# Restore the originally active sheet / selection on the Succeed sheet.
$succeed.Activate() | Out-Null
$succeed.Range("P17").Select() | Out-Null
